$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:H1, copying the formatting (style) of the existing headers
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# New boolean columns F:H (KNN/SVM/RF Outliers MAD flags)
$ws.Range("F2:H6").Value = $false
$ws.Range("G4").Value = $true

Write-Host "done"
